# Auto-generated edit script: updates Gilgamesh_Profits market-data cells
# across 8 item-category worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values are raw (non-formula) numbers scraped from a scheduled market-data
# job; this script simply overwrites / adds / removes the affected cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2447.1724
$ws.Range("I15").Value = 2447.1724
$ws.Range("K15").Value = 7341.5172
$ws.Range("M15").Value = -7172.5172
$ws.Range("H18").Value = 200000480
$ws.Range("I18").Value = 593
$ws.Range("K18").Value = 593
$ws.Range("M18").Value = -309
$ws.Range("H96").Value = 974.7143
$ws.Range("I96").Value = 868.5
$ws.Range("J96").Value = 1116.3334
$ws.Range("K96").Value = 2605.5
$ws.Range("L96").Value = 3349.0002
$ws.Range("M96").Value = -1232.5
$ws.Range("N96").Value = -6095.0002
$ws.Range("H98").Value = 4443.095
$ws.Range("I98").Value = 4166.8823
$ws.Range("K98").Value = 4166.8823
$ws.Range("M98").Value = -2668.8823
$ws.Range("H122").Value = 4443.095
$ws.Range("I122").Value = 4166.8823
$ws.Range("K122").Value = 12500.6469
$ws.Range("M122").Value = -10050.6469
$ws.Range("H127").Value = 532.375
$ws.Range("I127").Value = 328.33334
$ws.Range("J127").Value = 1144.5
$ws.Range("K127").Value = 985.0000200000001
$ws.Range("L127").Value = 3433.5
$ws.Range("M127").Value = 3974.99998
$ws.Range("N127").Value = -13353.5
$ws.Range("H131").Value = 564012.1
$ws.Range("I131").Value = 722765.4
$ws.Range("K131").Value = 2168296.2
$ws.Range("M131").Value = -2163256.2
$ws.Range("H132").Value = 5691.6294
$ws.Range("I132").Value = 5872.077
$ws.Range("K132").Value = 17616.231
$ws.Range("M132").Value = -15086.231
$ws.Range("H138").Value = 2788.9673
$ws.Range("I138").Value = 1317.8
$ws.Range("J138").Value = 2920.3215
$ws.Range("K138").Value = 3953.4
$ws.Range("L138").Value = 8760.9645
$ws.Range("M138").Value = 1186.6
$ws.Range("N138").Value = -19040.9645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2910.7715
$ws.Range("I32").Value = 2908.2058
$ws.Range("K32").Value = 2908.2058
$ws.Range("M32").Value = -2621.2058
$ws.Range("H45").Value = 107352.5
$ws.Range("I45").Value = 142136.67
$ws.Range("K45").Value = 142136.67
$ws.Range("M45").Value = -141759.67
$ws.Range("H61").Value = 2769.7646
$ws.Range("I61").Value = 1290.125
$ws.Range("K61").Value = 1290.125
$ws.Range("M61").Value = -1078.125
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws.Range("H122").Value = 2282.9312
$ws.Range("I122").Value = 2137.2222
$ws.Range("K122").Value = 6411.6666
$ws.Range("M122").Value = -3961.6666
$ws.Range("H136").Value = 2769.7646
$ws.Range("I136").Value = 1290.125
$ws.Range("K136").Value = 3870.375
$ws.Range("M136").Value = -1320.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2986
$ws.Range("I86").Value = 2592.7896
$ws.Range("J86").Value = 4231.1665
$ws.Range("K86").Value = 2592.7896
$ws.Range("L86").Value = 4231.1665
$ws.Range("M86").Value = -1469.7896
$ws.Range("N86").Value = -6477.1665
$ws.Range("H89").Value = 2986
$ws.Range("I89").Value = 2592.7896
$ws.Range("J89").Value = 4231.1665
$ws.Range("K89").Value = 12963.948
$ws.Range("L89").Value = 21155.8325
$ws.Range("M89").Value = -7347.948
$ws.Range("N89").Value = -32387.8325
$ws.Range("H99").Value = 5394
$ws.Range("I99").Value = 4972.8
$ws.Range("K99").Value = 4972.8
$ws.Range("M99").Value = -3474.8
$ws.Range("H134").Value = 1696.8846
$ws.Range("I134").Value = 1137.75
$ws.Range("K134").Value = 3413.25
$ws.Range("M134").Value = -878.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000462
$ws.Range("J6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("N6").Value = -2726
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40722
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42496
$ws.Range("H132").Value = 13338517
$ws.Range("I132").Value = 4336.231
$ws.Range("K132").Value = 13008.693
$ws.Range("M132").Value = -10478.693
$ws.Range("H134").Value = 4097.136
$ws.Range("I134").Value = 4191.737
$ws.Range("K134").Value = 12575.211
$ws.Range("M134").Value = -10040.211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 75.818184
$ws.Range("I4").Value = 76.38776
$ws.Range("K4").Value = 229.16328
$ws.Range("M4").Value = -117.16328
$ws.Range("H5").Value = 755.3333
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = ""
$ws.Range("H23").Value = 7152.3335
$ws.Range("I23").Value = 650
$ws.Range("J23").Value = 8452.799999999999
$ws.Range("K23").Value = 1950
$ws.Range("L23").Value = 25358.4
$ws.Range("M23").Value = -1715
$ws.Range("N23").Value = -25828.4
$ws.Range("H40").Value = 267.57144
$ws.Range("I40").Value = 255.6
$ws.Range("J40").Value = 297.5
$ws.Range("K40").Value = 1022.4
$ws.Range("L40").Value = 1190
$ws.Range("M40").Value = -953.4
$ws.Range("N40").Value = -1328
$ws.Range("H56").Value = 6590.6
$ws.Range("I56").Value = 6590.6
$ws.Range("K56").Value = 6590.6
$ws.Range("M56").Value = -6060.6
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H130").Value = 11912
$ws.Range("I130").Value = 11912
$ws.Range("K130").Value = 35736
$ws.Range("M130").Value = -30716
$ws.Range("H135").Value = 755.3333
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H138").Value = 3749.5
$ws.Range("I138").Value = 1500
$ws.Range("J138").Value = 5999
$ws.Range("K138").Value = 4500
$ws.Range("L138").Value = 17997
$ws.Range("M138").Value = 640
$ws.Range("N138").Value = -28277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2117.6428
$ws.Range("I132").Value = 1972.8462
$ws.Range("K132").Value = 5918.5386
$ws.Range("M132").Value = -3388.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2200.4
$ws.Range("I7").Value = 1875.5
$ws.Range("K7").Value = 1875.5
$ws.Range("M7").Value = -1763.5
$ws.Range("H100").Value = 4069.3
$ws.Range("I100").Value = 3159.6
$ws.Range("J100").Value = 4979
$ws.Range("K100").Value = 3159.6
$ws.Range("L100").Value = 4979
$ws.Range("M100").Value = -2618.6
$ws.Range("N100").Value = -6061
$ws.Range("H126").Value = 2200.4
$ws.Range("I126").Value = 1875.5
$ws.Range("K126").Value = 5626.5
$ws.Range("M126").Value = -3156.5
$ws.Range("H132").Value = 6731.364
$ws.Range("I132").Value = 4507.3335
$ws.Range("J132").Value = 9400.200000000001
$ws.Range("K132").Value = 13522.0005
$ws.Range("L132").Value = 28200.6
$ws.Range("M132").Value = -10992.0005
$ws.Range("N132").Value = -33260.60000000001
$ws.Range("H134").Value = 102388.5
$ws.Range("J134").Value = 102388.5
$ws.Range("L134").Value = 102388.5
$ws.Range("N134").Value = -112528.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4229
$ws.Range("I81").Value = 3520.6667
$ws.Range("K81").Value = 7041.3334
$ws.Range("M81").Value = -5980.3334
$ws.Range("H84").Value = 4229
$ws.Range("I84").Value = 3520.6667
$ws.Range("K84").Value = 35206.667
$ws.Range("M84").Value = -29902.667
$ws.Range("H122").Value = 6947157.5
$ws.Range("I122").Value = 2817.5
$ws.Range("K122").Value = 8452.5
$ws.Range("M122").Value = -6002.5
$ws.Range("H126").Value = 3199.4285
$ws.Range("I126").Value = 3149.3333
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 9447.999899999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -6977.999899999999
$ws.Range("N126").Value = -15440
$ws.Range("H141").Value = 68873.81
$ws.Range("J141").Value = 68798.8
$ws.Range("L141").Value = 68798.8
$ws.Range("N141").Value = -79158.8

